$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$xlPasteFormats = -4122

# Fill in the three new time-tracking entries (rows 24-26), matching the
# style/format already used by the existing rows above (Stunden / Datum /
# Zeitraum / Tasks).

# Row 24
$ws.Range("A24").Value = 1
$ws.Range("B24").Value = 43551
$ws.Range("C24").Value = "19:00:00-20:00"
$ws.Range("D24").Value = "Präsentation - block, delegate_to"

# Row 25
$ws.Range("A25").Value = 1.5
$ws.Range("B25").Value = 43552
$ws.Range("C25").Value = "10:00-11:30"
$ws.Range("D25").Value = "Präs. lookup, include_tasks, include_playbook"

# Row 26
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = 43552
$ws.Range("C26").Value = "16:16 - 17:15"
$ws.Range("D26").Value = "Präs. lookup, include_tasks, include_playbook"

# Copy number-format/alignment styling from the existing rows so that no new
# style definitions are introduced: B column uses the date format seen on
# B5:B23, C24/C26 use the time format seen on C13/C14/C15/C16/etc.
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B24:B26").PasteSpecial($xlPasteFormats)

$ws.Range("C13").Copy() | Out-Null
$ws.Range("C24").PasteSpecial($xlPasteFormats)
$ws.Range("C26").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# Update the view: scroll so row 10 is at top, and select D25 as active cell.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D25").Select()
